# Fully update for new ResolvePM poll
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# --- New raw poll input numbers (J3:O3) for the latest ResolvePM release ---
$ws.Range("J3").Value2 = 33
$ws.Range("K3").Value2 = 35
$ws.Range("L3").Value2 = 10
$ws.Range("M3").Value2 = 2
$ws.Range("N3").Value2 = 15
$ws.Range("O3").Value2 = 5

# Force a recalculation so P3 (the derived 2PP figure) is up to date before
# we copy it down as a plain value into the "latest ResolvePM" row.
$excel.Calculate()

# --- Shift the ResolvePM poll history down one slot ---
# Third ResolvePM <- Second ResolvePM (old row 10 data)
$ws.Range("B11").Value2 = $ws.Range("B10").Value2
$ws.Range("C11").Value2 = $ws.Range("C10").Value2
$ws.Range("D11").Value2 = $ws.Range("D10").Value2
$ws.Range("E11").Value2 = $ws.Range("E10").Value2
$ws.Range("H11").Value2 = $ws.Range("H10").Value2

# Second ResolvePM <- Latest ResolvePM (old row 9 data)
$ws.Range("B10").Value2 = $ws.Range("B9").Value2
$ws.Range("C10").Value2 = $ws.Range("C9").Value2
$ws.Range("D10").Value2 = $ws.Range("D9").Value2
$ws.Range("E10").Value2 = $ws.Range("E9").Value2
$ws.Range("H10").Value2 = $ws.Range("H9").Value2

# Latest ResolvePM <- brand new poll figures
$ws.Range("B9").Value2 = 53.704999999999998
$ws.Range("C9").Value2 = 52.332000000000001
$ws.Range("D9").Value2 = 58.515999999999998
$ws.Range("E9").Value2 = 48.938000000000002
$ws.Range("H9").Value2 = $ws.Range("P3").Value2

# --- Selection bookkeeping to match the saved workbook state ---
$ws.Range("E31").Select()

$excel.Calculate()
